$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.055887666666667
$ws.Range("H2").Value = 12.167663
$ws.Range("I2").Value = 0.4763357569530485
$ws.Range("J2").Value = 0.4763357569530485
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 102.9038579454919
$ws.Range("R2").Value = 926.1347215094271
$ws.Range("S2").Value = 0.2753640547400232
$ws.Range("T2").Value = 0.2753640547400232
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.055887666666667
$ws.Range("H3").Value = 12.167663
$ws.Range("I3").Value = 0.4763357569530485
$ws.Range("J3").Value = 0.4763357569530485
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 41.44992138796945
$ws.Range("R3").Value = 373.0492924917251
$ws.Range("S3").Value = 0.1109173032957848
$ws.Range("T3").Value = 0.1109173032957848
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.055887666666667
$ws.Range("H4").Value = 12.167663
$ws.Range("I4").Value = 0.4763357569530485
$ws.Range("J4").Value = 0.4763357569530485
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 33.65343048240445
$ws.Range("R4").Value = 302.88087434164
$ws.Range("S4").Value = 0.09005439891724056
$ws.Range("T4").Value = 0.09005439891724057
$ws.Range("I5").Value = 0.4564376967244237
$ws.Range("J5").Value = 0.4564376967244237
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 98.60523636760544
$ws.Range("R5").Value = 887.447127308449
$ws.Range("S5").Value = 0.2638612219880503
$ws.Range("T5").Value = 0.2638612219880503
$ws.Range("I6").Value = 0.4564376967244237
$ws.Range("J6").Value = 0.4564376967244237
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.1062839346075011
$ws.Range("T6").Value = 0.1062839346075011
$ws.Range("I7").Value = 0.4564376967244237
$ws.Range("J7").Value = 0.4564376967244237
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.08629254012887233
$ws.Range("T7").Value = 0.08629254012887234
$ws.Range("G8").Value = 0.5724183333333334
$ws.Range("I8").Value = 0.06722654632252778
$ws.Range("J8").Value = 0.06722654632252777
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 14.52309819693278
$ws.Range("R8").Value = 130.707883772395
$ws.Range("S8").Value = 0.03886286954385394
$ws.Range("T8").Value = 0.03886286954385394
$ws.Range("G9").Value = 0.5724183333333334
$ws.Range("I9").Value = 0.06722654632252778
$ws.Range("J9").Value = 0.06722654632252777
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("Q9").Value = 5.84993887101389
$ws.Range("S9").Value = 0.01565405728866775
$ws.Range("T9").Value = 0.01565405728866774
$ws.Range("G10").Value = 0.5724183333333334
$ws.Range("I10").Value = 0.06722654632252778
$ws.Range("J10").Value = 0.06722654632252777
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("Q10").Value = 4.749599143488889
$ws.Range("R10").Value = 42.7463922914
$ws.Range("S10").Value = 0.01270961949000609
$ws.Range("T10").Value = 0.01270961949000609
